# Fix bugs for other countries
#
# The 3rd tab ("Inflation rate (%)") actually contained data for the wrong
# indicator; it is repurposed into "Unemployment rate (%)" with fresh data.
# The 4th tab ("Current account balance (% of GDP)") is repurposed into the
# new "Inflation rate (%)" sheet, taking over the data that used to live on
# the old inflation-rate sheet.

$wb = $excel.ActiveWorkbook

$wsThird = $wb.Worksheets.Item("Inflation rate (%)")
$wsFourth = $wb.Worksheets.Item("Current account balance (% of G")

# Capture the old inflation-rate values (Value2 avoids COM Variant wrapper
# weirdness seen with plain .Value on read).
$oldB2 = $wsThird.Range("B2").Value2
$oldB3 = $wsThird.Range("B3").Value2
$oldB4 = $wsThird.Range("B4").Value2
$oldB5 = $wsThird.Range("B5").Value2
$oldB6 = $wsThird.Range("B6").Value2

# --- 3rd sheet: becomes "Unemployment rate (%)" with new data -------------
$wsThird.Range("B1").Value = "Unemployment rate (%)"
$wsThird.Range("B2").Value = 5.6
$wsThird.Range("B3").Value = 5.6
$wsThird.Range("B4").Value = 5.2
$wsThird.Range("B5").Value = 4.8
# Row 6 (2019) value is untouched by this fix.
$wsThird.Name = "Unemployment rate (%)"

# --- 4th sheet: becomes "Inflation rate (%)", reusing the old data --------
$wsFourth.Range("B1").Value = "Inflation rate (%)"
$wsFourth.Range("B2").Value = $oldB2
$wsFourth.Range("B3").Value = $oldB3
$wsFourth.Range("B4").Value = $oldB4
$wsFourth.Range("B5").Value = $oldB5
$wsFourth.Range("B6").Value = $oldB6
$wsFourth.Name = "Inflation rate (%)"
